$wb = $excel.ActiveWorkbook

# Update both sheets that contain this data table: "展览" and "全部类型"
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 27
    $ws.Range("F3").Value = 54
}
